$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values for column A (scenario names) used repeatedly across the new rows
$teacher      = "Login as a Teacher"
$student      = "Login as a Student"
$admin        = "Login as an Admin"
$teacherNeg1  = "Login as a Teacher (Negative-1)"
$teacherNeg2  = "Login as a Teacher (Negative-2)"
$studentNeg1  = "Login as a Student (Negative-1)"
$studentNeg2  = "Login as a Student (Negative-2)"
$adminNeg1    = "Login as an Admin (Negative-1)"
$adminNeg2    = "Login as an Admin (Negative-2)"

$passed = "PASSED"
$failed = "FAILED"
$chrome = "chrome"

# Rows 4-6: repeat of the positive scenarios (Teacher, Student, Admin) - PASSED
$rows = @(
    @(4,  $teacher,     $passed, $chrome),
    @(5,  $student,     $passed, $chrome),
    @(6,  $admin,       $passed, $chrome),
    @(7,  $teacherNeg1, $passed, $chrome),
    @(8,  $teacherNeg2, $passed, $chrome),
    @(9,  $studentNeg1, $failed, $chrome),
    @(10, $studentNeg2, $failed, $chrome),
    @(11, $adminNeg1,   $failed, $chrome),
    @(12, $adminNeg2,   $failed, $chrome),
    @(13, $teacher,     $passed, $chrome),
    @(14, $student,     $passed, $chrome),
    @(15, $admin,       $passed, $chrome),
    @(16, $teacherNeg1, $passed, $chrome),
    @(17, $teacherNeg2, $passed, $chrome),
    @(18, $studentNeg1, $passed, $chrome),
    @(19, $studentNeg2, $passed, $chrome),
    @(20, $adminNeg1,   $passed, $chrome),
    @(21, $adminNeg2,   $passed, $chrome)
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 1).Value = $r[1]
    $ws.Cells.Item($rowNum, 2).Value = $r[2]
    $ws.Cells.Item($rowNum, 3).Value = $r[3]
}
